# Update the "Date" column (column B) of the RAD test-data sheet with the
# timestamps captured from the latest Katalon test run (Mon Sep 11 2023,
# 13:48:47 EDT - 13:53:34 EDT), replacing the previous run's timestamps
# (Fri Sep 08 2023, 17:54:30 EDT - 17:59:45 EDT) that were recorded in
# B2:B29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamps = @(
    "Mon Sep 11 13:48:47 EDT 2023",
    "Mon Sep 11 13:48:58 EDT 2023",
    "Mon Sep 11 13:49:08 EDT 2023",
    "Mon Sep 11 13:49:19 EDT 2023",
    "Mon Sep 11 13:49:29 EDT 2023",
    "Mon Sep 11 13:49:40 EDT 2023",
    "Mon Sep 11 13:49:51 EDT 2023",
    "Mon Sep 11 13:50:02 EDT 2023",
    "Mon Sep 11 13:50:13 EDT 2023",
    "Mon Sep 11 13:50:23 EDT 2023",
    "Mon Sep 11 13:50:34 EDT 2023",
    "Mon Sep 11 13:50:44 EDT 2023",
    "Mon Sep 11 13:50:55 EDT 2023",
    "Mon Sep 11 13:51:06 EDT 2023",
    "Mon Sep 11 13:51:16 EDT 2023",
    "Mon Sep 11 13:51:27 EDT 2023",
    "Mon Sep 11 13:51:38 EDT 2023",
    "Mon Sep 11 13:51:48 EDT 2023",
    "Mon Sep 11 13:51:59 EDT 2023",
    "Mon Sep 11 13:52:10 EDT 2023",
    "Mon Sep 11 13:52:20 EDT 2023",
    "Mon Sep 11 13:52:30 EDT 2023",
    "Mon Sep 11 13:52:41 EDT 2023",
    "Mon Sep 11 13:52:51 EDT 2023",
    "Mon Sep 11 13:53:02 EDT 2023",
    "Mon Sep 11 13:53:13 EDT 2023",
    "Mon Sep 11 13:53:23 EDT 2023",
    "Mon Sep 11 13:53:34 EDT 2023"
)

$startRow = 2
for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newTimestamps[$i]
}
